$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated ROI figures after adding article calculations.
# Red highlight (Interior.Color = 255) marks the "good ROI" cells that
# reuse the workbook's existing red-fill style; Style = "Normal" clears
# that highlight back to the default (no-fill) style.

# Row 2 - Argentina
$ws.Range("D2").Value = 275
$ws.Range("D2").Interior.Color = 255
$ws.Range("G2").Value = 53
$ws.Range("G2").Interior.Color = 255
$ws.Range("I2").Value = -43.75
$ws.Range("K2").Value = 4.57
$ws.Range("M2").Value = -0.27
$ws.Range("O2").Value = -2.73
$ws.Range("Q2").Value = 26.14
$ws.Range("T2").Value = 333
$ws.Range("T2").Interior.Color = 255
$ws.Range("W2").Value = 42
$ws.Range("X2").Value = 13.8
$ws.Range("AA2").Value = -0.87
$ws.Range("AB2").Value = 15

# Row 7 - Denmark
$ws.Range("B7").Value = 8.76
$ws.Range("E7").Value = 12.32
$ws.Range("H7").Value = -0.67
$ws.Range("K7").Value = 12.36
$ws.Range("M7").Value = -17.06
$ws.Range("O7").Value = -7.41
$ws.Range("Q7").Value = 12.2
$ws.Range("S7").Value = 29.21
$ws.Range("U7").Value = 3.33
$ws.Range("X7").Value = -11.74
$ws.Range("AA7").Value = 17.38
$ws.Range("AB7").Value = 45

# Row 21 - Mexico2
$ws.Range("B21").Value = 25.17
$ws.Range("D21").Value = -32.75
$ws.Range("E21").Value = 12.73
$ws.Range("G21").Value = -53.1
$ws.Range("H21").Value = 1.59
$ws.Range("J21").Value = 24.45
$ws.Range("K21").Value = -41
$ws.Range("M21").Value = -14.89
$ws.Range("O21").Value = -6.97
$ws.Range("P21").Value = 9.21
$ws.Range("Q21").Value = -44.14
$ws.Range("R21").Value = -59.11
$ws.Range("S21").Value = -22.83
$ws.Range("U21").Value = -3.32
$ws.Range("W21").Value = -39.27
$ws.Range("X21").Value = 9.92
$ws.Range("Z21").Value = 8.5
$ws.Range("AA21").Value = -41.83
$ws.Range("AB21").Value = 40

# Row 24 - Poland
$ws.Range("B24").Value = -18.85
$ws.Range("D24").Value = 23.24
$ws.Range("E24").Value = -6.85
$ws.Range("F24").Value = 1.36
$ws.Range("F24").Interior.Color = 255
$ws.Range("H24").Value = -10.09
$ws.Range("J24").Value = -21.23
$ws.Range("K24").Value = -2.17
$ws.Range("L24").Value = -4.17
$ws.Range("M24").Value = -0.58
$ws.Range("O24").Value = -2.52
$ws.Range("P24").Value = -22.81
$ws.Range("R24").Value = -2.17
$ws.Range("R24").Style = "Normal"
$ws.Range("S24").Value = -36.49
$ws.Range("T24").Value = -12.43
$ws.Range("U24").Value = -13.29
$ws.Range("W24").Value = -9.84
$ws.Range("X24").Value = 6.44
$ws.Range("Z24").Value = 12.59
$ws.Range("AA24").Value = -5.16
$ws.Range("AB24").Value = 115

# Row 26 - Portugal2
$ws.Range("B26").Value = -0.02
$ws.Range("B26").Style = "Normal"
$ws.Range("E26").Value = -2.17
$ws.Range("H26").Value = 0.28
$ws.Range("J26").Value = -29.07
$ws.Range("M26").Value = 0.24
$ws.Range("M26").Interior.Color = 255
$ws.Range("O26").Value = -2.96
$ws.Range("P26").Value = -17.72
$ws.Range("R26").Value = 10.37
$ws.Range("U26").Value = -0.11
$ws.Range("X26").Value = 0.72
$ws.Range("Z26").Value = 22.74
$ws.Range("AB26").Value = 120

# Row 33 - Spain2
$ws.Range("B33").Value = -1.11
$ws.Range("E33").Value = 0.76
$ws.Range("H33").Value = -6.67
$ws.Range("K33").Value = -8.24
$ws.Range("M33").Value = -8.31
$ws.Range("O33").Value = -2.56
$ws.Range("Q33").Value = -8.44
$ws.Range("S33").Value = -14.77
$ws.Range("U33").Value = -1.35
$ws.Range("X33").Value = -8.15
$ws.Range("AA33").Value = -2.47
$ws.Range("AB33").Value = 163

# Row 36 - Turkey
$ws.Range("B36").Value = 8.04
$ws.Range("D36").Value = -30.25
$ws.Range("E36").Value = -1.82
$ws.Range("G36").Value = -20.17
$ws.Range("H36").Value = -5.9
$ws.Range("J36").Value = -13.65
$ws.Range("K36").Value = -27.7
$ws.Range("M36").Value = -9.39
$ws.Range("O36").Value = -9.3
$ws.Range("P36").Value = -15.45
$ws.Range("Q36").Value = -26.68
$ws.Range("R36").Value = -5.62
$ws.Range("S36").Value = -14.76
$ws.Range("U36").Value = 0.45
$ws.Range("W36").Value = -9.17
$ws.Range("X36").Value = -1.79
$ws.Range("AA36").Value = -6.34
$ws.Range("AB36").Value = 144
